# Updated data following Rachel's reporting.
#
# The "crime" sheet contained several duplicate / superseded rows for the
# same recipient+office (older project descriptions/amounts that were
# later corrected or retracted). Remove those rows outright; the "courts"
# sheet is untouched (its cells keep their values - only the shared
# string table shrinks as a side effect of saving).
#
# Rows removed (1-based worksheet row numbers, header = row 1):
#   5  - Cleveland Police Department / $355,400.00 (RTCC description)
#   9  - Cuyahoga County Prosecutor's Office / $250,670.19 (NIBIN descr.)
#   11 - Cuyahoga County Sheriff's Office / $163,000.00 (hot spot descr.)
#   13 - East Cleveland Police Department row
#   14 - Euclid Police Department / $174,771.75 (GVRT description)
#   22 - Westlake Police Department row
#
# Deleted from the bottom up so earlier row numbers stay valid while the
# later deletes are still pending.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crime")

$ws.Rows(22).Delete()
$ws.Rows(14).Delete()
$ws.Rows(13).Delete()
$ws.Rows(11).Delete()
$ws.Rows(9).Delete()
$ws.Rows(5).Delete()
